# Reorder the block-order header columns and corresponding one-hot data
# to match the new canonical column ordering:
#   living_rooms_1, living_rooms_2, bedrooms_1, bedrooms_2, kitchens_1, kitchens_2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "kitchens_1"
$ws.Range("F1").Value = "kitchens_2"

# --- Data rows (rows 2-7), one-hot encoded values ---
$data = @(
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
